# Add a new "LEADING CAUSES TREND" row to the TRENDS group of the tab
# parameter-linkage table (queryParameter_Linkage.xlsx / sheet "tab"),
# inserted right before the existing "DISPARITIES" row (i.e. becomes the
# new row 15, pushing everything below it down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tab")

# Insert a new blank row at position 15 (shifts old row 15.."DISPARITIES"
# and everything after it down to 16..25).
$ws.Rows.Item(15).Insert()

# Populate the new row with the new "topTrendsTab" / leading causes trend
# tab entry, matching the shape of the other TRENDS rows (10-14).
$ws.Range("A15").Value = "TRENDS"
$ws.Range("B15").Value = "LEADING CAUSES TREND"
$ws.Range("C15").Value = "trends"
$ws.Range("D15").Value = "topTrendsTab"
$ws.Range("E15").Value = "leadingcausestrend"

# Match the author's final cursor/selection position.
[void]$ws.Range("I15").Select()
